# Update "想去人数" (F) / "最低票价" (G) figures to the latest scrape snapshot.
$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 155
$ws.Range("F3").Value = 168
$ws.Range("F4").Value = 2192
$ws.Range("F5").Value = 4337
$ws.Range("F6").Value = 577
$ws.Range("F7").Value = 1067
$ws.Range("F8").Value = 1346
$ws.Range("F9").Value = 668
$ws.Range("F10").Value = 381
$ws.Range("F11").Value = 106
$ws.Range("F12").Value = 407
$ws.Range("F13").Value = 674401
$ws.Range("F14").Value = 1673
$ws.Range("F15").Value = 579
$ws.Range("F16").Value = 1491
$ws.Range("F17").Value = 677
$ws.Range("F18").Value = 548
$ws.Range("F19").Value = 1306
$ws.Range("F20").Value = 2304
$ws.Range("F21").Value = 1163
$ws.Range("F22").Value = 2726
$ws.Range("F23").Value = 1582
$ws.Range("F24").Value = 880
$ws.Range("F25").Value = 1575
$ws.Range("F26").Value = 541
$ws.Range("F27").Value = 1100
$ws.Range("F28").Value = 910
$ws.Range("F29").Value = 1104
$ws.Range("F31").Value = 91
$ws.Range("F32").Value = 2051
$ws.Range("F34").Value = 1341
$ws.Range("F35").Value = 3039
$ws.Range("F38").Value = 49
$ws.Range("F39").Value = 203
$ws.Range("F40").Value = 2628
$ws.Range("F41").Value = 216
$ws.Range("F43").Value = 3178
$ws.Range("F47").Value = 164
$ws.Range("F49").Value = 24

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F10").Value = 496
$ws.Range("F11").Value = 144991
$ws.Range("F12").Value = 144992
$ws.Range("F19").Value = 342
$ws.Range("F22").Value = 177
$ws.Range("F23").Value = 86
$ws.Range("F26").Value = 639
$ws.Range("F31").Value = 373
$ws.Range("G31").Value = 180
$ws.Range("F32").Value = 281
$ws.Range("F33").Value = 13
$ws.Range("F34").Value = 73
$ws.Range("F35").Value = 73
$ws.Range("F41").Value = 189
$ws.Range("F42").Value = 90

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 3138
$ws.Range("F7").Value = 836
$ws.Range("F8").Value = 1224
$ws.Range("F9").Value = 647
$ws.Range("F10").Value = 1617
$ws.Range("F11").Value = 141
$ws.Range("F12").Value = 2026

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 836
$ws.Range("F3").Value = 647
$ws.Range("F4").Value = 155
$ws.Range("F5").Value = 1617
$ws.Range("F6").Value = 168
$ws.Range("F7").Value = 2192
$ws.Range("F8").Value = 141
$ws.Range("F9").Value = 2026
$ws.Range("F10").Value = 4337
$ws.Range("F11").Value = 577
$ws.Range("F12").Value = 1346
$ws.Range("F13").Value = 668
$ws.Range("F14").Value = 381
$ws.Range("F15").Value = 407
$ws.Range("F16").Value = 674411
$ws.Range("F18").Value = 496
$ws.Range("F19").Value = 1673
$ws.Range("F20").Value = 144992
$ws.Range("F21").Value = 1491
$ws.Range("F22").Value = 677
$ws.Range("F23").Value = 548
$ws.Range("F24").Value = 1306
$ws.Range("F25").Value = 2304
$ws.Range("F26").Value = 1163
$ws.Range("F27").Value = 2726
$ws.Range("F28").Value = 1582
$ws.Range("F29").Value = 880
$ws.Range("F31").Value = 1575
$ws.Range("F32").Value = 541
$ws.Range("F33").Value = 177
$ws.Range("F34").Value = 1100
$ws.Range("F35").Value = 910
$ws.Range("F36").Value = 1104
$ws.Range("F37").Value = 91
$ws.Range("F38").Value = 2051
$ws.Range("F39").Value = 1341
$ws.Range("F40").Value = 3039
$ws.Range("F43").Value = 373
$ws.Range("G43").Value = 180
$ws.Range("F44").Value = 281
$ws.Range("F45").Value = 73
$ws.Range("F46").Value = 2628
$ws.Range("F49").Value = 3178
$ws.Range("F52").Value = 164
$ws.Range("F54").Value = 24
